$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'65.172.98"
$ws.Range("E2").Value = "  +1.57%  "
$ws.Range("D3").Value = "'3.134.96"
$ws.Range("E3").Value = "  +2.16%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").Value = "'570.52"
$ws.Range("E5").Value = "  +2.04%  "
$ws.Range("D6").Value = "'149.38"
$ws.Range("E6").Value = "  +2.70%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'3.131.92"
$ws.Range("E8").Value = "  +2.12%  "
$ws.Range("D9").Value = "'0.527"
$ws.Range("E9").Value = "  +5.07%  "
$ws.Range("E10").Value = "  +6.88%  "
$ws.Range("D11").Value = "'6.17"
$ws.Range("E11").Value = "  +0.28%  "
$ws.Range("D12").Value = "'0.500"
$ws.Range("E12").Value = "  +6.75%  "
$ws.Range("D13").Value = "'0.0000259"
$ws.Range("E13").Value = "  +13.92%  "
$ws.Range("D14").Value = "'37.13"
$ws.Range("E14").Value = "  +5.53%  "
$ws.Range("D15").Value = "'3.651.03"
$ws.Range("E15").Value = "  +2.30%  "
$ws.Range("D16").Value = "'65.224.98"
$ws.Range("E16").Value = "  +1.70%  "
$ws.Range("E17").Value = "  +5.57%  "
$ws.Range("D18").Value = "'3.140.87"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D20").Value = "'511.25"
$ws.Range("E20").Value = "  +6.86%  "
$ws.Range("E21").Value = "  +6.64%  "
$ws.Range("D22").Value = "'15.50"
$ws.Range("E22").Value = "  +14.25%  "
$ws.Range("D23").Value = "'0.724"
$ws.Range("E23").Value = "  +7.08%  "
$ws.Range("D24").Value = "'7.80"
$ws.Range("E24").Value = "  +3.12%  "
$ws.Range("D25").Value = "'85.45"
$ws.Range("E25").Value = "  +4.55%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E27").Value = "  +3.88%  "
$ws.Range("D28").Value = "'8.74"
$ws.Range("E28").Value = "  +7.77%  "
$ws.Range("E29").Value = "  +3.65%  "
$ws.Range("D30").Value = "'27.97"
$ws.Range("E30").Value = "  +6.74%  "
$ws.Range("D31").Value = "'1.00"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("B32").Value = "Mantle"
$ws.Range("C32").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D32").Value = "'1.19"
$ws.Range("E32").Value = "  +3.60%  "
$ws.Range("B33").Value = "Stacks"
$ws.Range("C33").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D33").Value = "'2.66"
$ws.Range("E33").Value = "  +7.02%  "
$ws.Range("E34").Value = "  +8.24%  "
$ws.Range("D35").Value = "'6.59"
$ws.Range("E35").Value = "  +6.67%  "
$ws.Range("D36").Value = "'55.69"
$ws.Range("E36").Value = "  +1.40%  "
$ws.Range("D37").Value = "'473.96"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("E38").Value = "  +4.27%  "
$ws.Range("D39").Value = "'0.0853"
$ws.Range("E39").Value = "  +2.89%  "
$ws.Range("D40").Value = "'2.98"
$ws.Range("E40").Value = "  -1.56%  "
$ws.Range("D41").Value = "'3.126.71"
$ws.Range("E41").Value = "  +5.26%  "
$ws.Range("B42").Value = "Cosmos"
$ws.Range("C42").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D42").Value = "'8.61"
$ws.Range("E42").Value = "  +4.15%  "
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").Value = "'0.121"
$ws.Range("E43").Value = "  +5.33%  "
$ws.Range("E44").Value = "  +9.62%  "
$ws.Range("E45").Value = "  +13.06%  "
$ws.Range("D46").Value = "'29.09"
$ws.Range("E46").Value = "  +4.27%  "
$ws.Range("E47").Value = "  -0.05%  "
$ws.Range("B48").Value = "PEPE"
$ws.Range("C48").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D48").Value = "'0.0₃0561"
$ws.Range("E48").Value = "  +8.63%  "
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").Value = "'0.116"
$ws.Range("E49").Value = "  +3.59%  "
$ws.Range("D50").Value = "'2.29"
$ws.Range("E50").Value = "  +10.13%  "
$ws.Range("D51").Value = "'117.95"
$ws.Range("E51").Value = "  -1.66%  "
